# [Fix]: Archivos QA gateway, Archivos Excel Catalogo general
# Rename the "Reactivos" (reagents) sheet/template into a generic "Catálogos"
# (catalogs) template: drop the "Nombre / Clave Contaq / Nombre Contaq" columns,
# add "Largo" / "Ancho" columns, keep "Clave" and "Activo", and drop the
# now-unused 5th (E) column altogether.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet.
$ws.Name = "Catálogos"

# Rename the sheet-scoped defined name and shrink its range from E to D.
$n = $wb.Names.Item(1)
$n.Name = "Catalogos"
$n.RefersTo = "=Catálogos!`$A`$4:`$D`$5"

# Header row (row 3): keep Clave/Activo, swap in Largo/Ancho, drop column E.
$ws.Range("B3").Value = "Largo"
$ws.Range("C3").Value = "Ancho"
$ws.Range("D3").Value = "Activo"
$ws.Range("E3").Clear()

# Placeholder row (row 4): keep item.Clave/item.Activo, swap placeholders,
# drop column E.
$ws.Range("B4").Value = "{{item.Largo}}"
$ws.Range("C4").Value = "{{item.Ancho}}"
$ws.Range("D4").Value = "{{item.Activo}}"
$ws.Range("E4").Clear()

# Restore the selection to the cell it was left on when the file was saved.
$ws.Range("G1").Select()
